$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 742 (shifts existing rows 742:783 down to 743:784)
$ws.Rows.Item(742).Insert()

# Populate the newly inserted row with the new log entry.
# Force column A to stay literal text ("2026/02/02" must not be
# auto-converted into a date serial number), then drop the temporary
# number-format override so the cell keeps the sheet's plain/no-style look.
$ws.Cells.Item(742, 1).NumberFormat = "@"
$ws.Cells.Item(742, 1).Value = "2026/02/02"
$ws.Cells.Item(742, 1).ClearFormats()

$ws.Cells.Item(742, 2).Value = "月"
$ws.Cells.Item(742, 3).Value = 13
$ws.Cells.Item(742, 4).Value = 24
